# Update cryptocurrency price/volume figures per the latest GitHub Actions data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.815.08"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +0.20%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.646.50"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.23%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +1.04%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'217.43"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +0.67%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("E6").Value = "'  -0.49%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  +0.95%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  -0.85%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  -0.59%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'19.20"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -1.03%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.0844"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -0.25%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'1.871.41"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -0.50%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'1.638.24"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -0.73%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'4.18"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -0.95%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("E15").Value = "'  -1.42%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'64.69"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -2.91%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'26.819.27"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -0.01%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("E18").Value = "'  -2.37%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'214.29"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -3.06%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'  +0.92%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'4.35"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -1.31%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'  +12.81%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'6.29"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -1.58%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'9.36"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -2.37%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'145.13"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -1.62%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'1.02"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +1.24%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  -2.68%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'7.10"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -0.01%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'15.68"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -1.92%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  -1.84%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  +0.39%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  -3.55%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  -2.96%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'1.283.71"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -1.17%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "'  -2.95%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  +1.69%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D38").Value = "'0.537"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +1.80%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.827"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -0.86%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  +0.94%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.815"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -0.18%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'2.23"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -0.94%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'5.35"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -1.69%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'1.797.17"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +0.24%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'91.37"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -2.81%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'59.66"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -2.42%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'  -1.03%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.0₆0103"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -1.80%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.0519"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +0.24%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'7.65"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -1.92%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.0976"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -0.53%  "
$ws.Range("E51").Style = "Normal"
